$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 300
$ws.Range("I18").Value = 300
$ws.Range("K18").Value = 300
$ws.Range("M18").Value = -16

$ws.Range("H28").Value = 489.6
$ws.Range("J28").Value = 1064.5
$ws.Range("L28").Value = 1064.5
$ws.Range("N28").Value = -2034.5

$ws.Range("H34").Value = 3185.7144
$ws.Range("I34").Value = 3050
$ws.Range("K34").Value = 3050
$ws.Range("M34").Value = -2847

$ws.Range("H36").Value = 3185.7144
$ws.Range("I36").Value = 3050
$ws.Range("K36").Value = 3050
$ws.Range("M36").Value = -2335

$ws.Range("H43").Value = 6633.2
$ws.Range("J43").Value = 7307.6924
$ws.Range("L43").Value = 7307.6924
$ws.Range("N43").Value = -7445.6924

$ws.Range("H55").Value = 416.63635
$ws.Range("J55").Value = 261.25
$ws.Range("L55").Value = 261.25
$ws.Range("N55").Value = -689.25

$ws.Range("H70").Value = 669.41174
$ws.Range("I70").Value = 590
$ws.Range("J70").Value = 782.8570999999999
$ws.Range("K70").Value = 1770
$ws.Range("L70").Value = 2348.5713
$ws.Range("M70").Value = -1500
$ws.Range("N70").Value = -2888.5713

$ws.Range("H73").Value = 669.41174
$ws.Range("I73").Value = 590
$ws.Range("J73").Value = 782.8570999999999
$ws.Range("K73").Value = 1770
$ws.Range("L73").Value = 2348.5713
$ws.Range("M73").Value = -834
$ws.Range("N73").Value = -4220.5713

$ws.Range("H80").Value = 2035.25
$ws.Range("I80").Value = 2280.3333
$ws.Range("J80").Value = 1300
$ws.Range("K80").Value = 6840.999899999999
$ws.Range("L80").Value = 3900
$ws.Range("M80").Value = -5842.999899999999
$ws.Range("N80").Value = -5896

$ws.Range("H83").Value = 2035.25
$ws.Range("I83").Value = 2280.3333
$ws.Range("J83").Value = 1300
$ws.Range("K83").Value = 20522.9997
$ws.Range("L83").Value = 11700
$ws.Range("M83").Value = -15530.9997
$ws.Range("N83").Value = -21684

$ws.Range("H86").Value = 11205.3
$ws.Range("I86").Value = 12729.8
$ws.Range("K86").Value = 12729.8
$ws.Range("M86").Value = -11606.8

$ws.Range("H89").Value = 11205.3
$ws.Range("I89").Value = 12729.8
$ws.Range("K89").Value = 63649
$ws.Range("M89").Value = -58033

$ws.Range("H98").Value = 1588.3871
$ws.Range("I98").Value = 1704.8928
$ws.Range("K98").Value = 1704.8928
$ws.Range("M98").Value = -206.8928000000001

$ws.Range("H112").Value = 5663.727
$ws.Range("J112").Value = 2891.7778
$ws.Range("L112").Value = 8675.3334
$ws.Range("N112").Value = -10891.3334

$ws.Range("H122").Value = 1588.3871
$ws.Range("I122").Value = 1704.8928
$ws.Range("K122").Value = 5114.678400000001
$ws.Range("M122").Value = -2664.678400000001

$ws.Range("H137").Value = 5765.737
$ws.Range("I137").Value = 2163.348
$ws.Range("K137").Value = 6490.044
$ws.Range("M137").Value = -3940.044

$ws.Range("H138").Value = 2628.4783
$ws.Range("I138").Value = 1977.0555
$ws.Range("J138").Value = 3047.25
$ws.Range("K138").Value = 5931.166499999999
$ws.Range("L138").Value = 9141.75
$ws.Range("M138").Value = -791.1664999999994
$ws.Range("N138").Value = -19421.75

$ws.Range("H141").Value = 3269.25
$ws.Range("I141").Value = 3269.25
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 9807.75
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -4627.75
$ws.Range("N141").ClearContents()


# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 10541.286
$ws.Range("I31").Value = 13003.5
$ws.Range("J31").Value = 7258.3335
$ws.Range("K31").Value = 13003.5
$ws.Range("L31").Value = 7258.3335
$ws.Range("M31").Value = -12709.5
$ws.Range("N31").Value = -7846.3335

$ws.Range("H32").Value = 6761064
$ws.Range("I32").Value = 10872561
$ws.Range("J32").Value = 6461.4287
$ws.Range("K32").Value = 10872561
$ws.Range("L32").Value = 6461.4287
$ws.Range("M32").Value = -10872274
$ws.Range("N32").Value = -7035.4287

$ws.Range("H63").Value = 3501
$ws.Range("I63").Value = 3501
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 3501
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -2815
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 3501
$ws.Range("I66").Value = 3501
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 17505
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -14073
$ws.Range("N66").ClearContents()

$ws.Range("H74").Value = 20841.227
$ws.Range("I74").Value = 1845.8889
$ws.Range("J74").Value = 47142.46
$ws.Range("K74").Value = 1845.8889
$ws.Range("L74").Value = 47142.46
$ws.Range("M74").Value = -971.8888999999999
$ws.Range("N74").Value = -48890.46

$ws.Range("H77").Value = 20841.227
$ws.Range("I77").Value = 1845.8889
$ws.Range("J77").Value = 47142.46
$ws.Range("K77").Value = 9229.4445
$ws.Range("L77").Value = 235712.3
$ws.Range("M77").Value = -4861.4445
$ws.Range("N77").Value = -244448.3

$ws.Range("H88").Value = 2041.75
$ws.Range("J88").Value = 2055
$ws.Range("L88").Value = 2055
$ws.Range("N88").Value = -2867

$ws.Range("H91").Value = 2041.75
$ws.Range("J91").Value = 2055
$ws.Range("L91").Value = 2055
$ws.Range("N91").Value = -4863

$ws.Range("H105").Value = 34500
$ws.Range("J105").Value = 34500
$ws.Range("L105").Value = 34500
$ws.Range("N105").Value = -41488


# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 884.0476
$ws.Range("I86").Value = 878.75
$ws.Range("K86").Value = 878.75
$ws.Range("M86").Value = 244.25

$ws.Range("H89").Value = 884.0476
$ws.Range("I89").Value = 878.75
$ws.Range("K89").Value = 4393.75
$ws.Range("M89").Value = 1222.25

$ws.Range("H94").Value = 1275.16
$ws.Range("I94").Value = 1277.6316
$ws.Range("K94").Value = 1277.6316
$ws.Range("M94").Value = -826.6315999999999

$ws.Range("H99").Value = 10765.861
$ws.Range("I99").Value = 11229.454
$ws.Range("K99").Value = 11229.454
$ws.Range("M99").Value = -9731.454

$ws.Range("H102").Value = 14221.571
$ws.Range("I102").Value = 14221.571
$ws.Range("K102").Value = 14221.571
$ws.Range("M102").Value = -10976.571


# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H93").Value = 20000
$ws.Range("I93").Value = 20000
$ws.Range("K93").Value = 20000
$ws.Range("M93").Value = -18128


# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 590.3077
$ws.Range("I40").Value = 605.0909
$ws.Range("J40").Value = 509
$ws.Range("K40").Value = 2420.3636
$ws.Range("L40").Value = 2036
$ws.Range("M40").Value = -2351.3636
$ws.Range("N40").Value = -2174

$ws.Range("H112").Value = 6869
$ws.Range("J112").Value = 10725
$ws.Range("L112").Value = 32175
$ws.Range("N112").Value = -34391

$ws.Range("H131").Value = 1478.24
$ws.Range("J131").Value = 1480.75
$ws.Range("L131").Value = 4442.25
$ws.Range("N131").Value = -14522.25


# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4162
$ws.Range("I113").Value = 3296.3333
$ws.Range("K113").Value = 3296.3333
$ws.Range("M113").Value = -1126.3333

$ws.Range("H132").Value = 815359
$ws.Range("J132").Value = 2437437.2
$ws.Range("L132").Value = 7312311.600000001
$ws.Range("N132").Value = -7317371.600000001


# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 10903.23
$ws.Range("I68").Value = 10057.75
$ws.Range("K68").Value = 10057.75
$ws.Range("M68").Value = -9308.75

$ws.Range("H71").Value = 10903.23
$ws.Range("I71").Value = 10057.75
$ws.Range("K71").Value = 50288.75
$ws.Range("M71").Value = -46544.75

$ws.Range("H93").Value = 3858.9033
$ws.Range("I93").Value = 2459
$ws.Range("J93").Value = 11138.4
$ws.Range("K93").Value = 2459
$ws.Range("L93").Value = 11138.4
$ws.Range("M93").Value = -1211
$ws.Range("N93").Value = -13634.4


# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 17433.182
$ws.Range("J62").Value = 8616.333000000001
$ws.Range("L62").Value = 8616.333000000001
$ws.Range("N62").Value = -9864.333000000001

$ws.Range("H65").Value = 17433.182
$ws.Range("J65").Value = 8616.333000000001
$ws.Range("L65").Value = 43081.665
$ws.Range("N65").Value = -49321.665

$ws.Range("H132").Value = 519656
$ws.Range("I132").Value = 2250.647
$ws.Range("J132").Value = 2718628.8
$ws.Range("K132").Value = 6751.941
$ws.Range("L132").Value = 8155886.399999999
$ws.Range("M132").Value = -4221.941
$ws.Range("N132").Value = -8160946.399999999

